$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Civ")

# Add a "Description" column (C) next to Civilization/Starting Tech with a blurb for each civ
$ws.Range('C1').Value = 'Description'
$ws.Range("C1").Font.Bold = $true

$ws.Range('C2').Value = 'The Americans start with a free random great person in their capital. Each time the Americans convert 3 trade into production, they receive 2 production instead of 1.'  # Americans
$ws.Range('C3').Value = 'The Arabs start the game with one of each resource from the market. For each resource token the Arabs spend for any reason, they gain 1 culture. Each time the Arabs invest a coin, they advance once space up the culture track for free'  # Arabs
$ws.Range('C4').Value = 'Each time the Aztecs win a battle, they gain 3 trade. After the Aztecs fight a battle, they gain 1 culture for each unit that was killed in the battle. Each time the Aztecs gain a great person, they may build 2 unlocked units of their choice for free.'  # Aztecs
$ws.Range('C5').Value = 'The Chinese start with city walls in their capital. The Chinese gain 3 culture each time they explore a hut or conquer a village. The Chinese may save one of their killed unit after each battle, returning it to their standing force'  # Chinese
$ws.Range('C6').Value = 'The Egyptians start with a free random ancient wonder in their capital. Once per turn, during City Management, the Egyptians may build an unlocked building for free by using an action.'  # Egyptians
$ws.Range('C7').Value = 'English armies may gather icons and resources as though they where scouts. English figures may cross water, but may not end their movement in it.'  # England
$ws.Range('C8').Value = 'France start with +2 in combat bonus. The maximum number of social policies the French can adopt is increased by 1. The French starts with 1 extra social policy.'  # French
$ws.Range('C9').Value = 'The Germans start with 2 extra infantry units. After setup, each time the Germans research a tech that upgrades or unlocked a unit, they build one of that unit for free and gain one resource of their choice from the market.'  # Germans
$ws.Range('C10').Value = 'Each time the Greeks gain a great person, they draw one extra great person, keeping one and discarding the other.'  # Greeks
$ws.Range('C11').Value = 'The Indians start with a metropolis as their capital. When the Indians spend a resource, they may use it as incense, silk, iron or wheat. When the Indians devote a city to the arts, it produces 1 extra culture for each square containing a resource (silk, iron etc) in its outskirts.'  # Indians
$ws.Range('C12').Value = 'Japanese infantry units have +1 in strength. The Japanese require 3 less trade to research new techs of any level. '  # Japanese
$ws.Range('C13').Value = 'The Mongols start with 2 extra mounted units. When attacking, if the Mongols gain loot, they gain 1 extra loot.'  # Mongols
$ws.Range('C14').Value = 'The Romans advance one space on the culture track for free each time they build a wonder or city, and each time they conquer a city or village.'  # Romans
$ws.Range('C15').Value = 'The Russians start with 2 extra armies, and their stacking limit is increased by 1. Once per turn, the Russians maybe move an army or scout into an enemy city and sacrifice that figure to research a tech known by that civilization for free. Armies sacrificed this way cannot be also attack.'  # Russians
$ws.Range('C16').Value = 'The Spanish start with 2 scouts. The travel speed of the Spanish is increased by 1. When the Spanish discover an unexplored map tile, they may immediately build a basic (non-upgraded) building in any of their cities for free, even if they haven''t unlocked that building.'  # Spanish
$ws.Range('C17').Value = 'The Zulu start with 2 extra artillery units. Zulu armies explore barbarian villages without a battle. The Zulu may build cities adjacent to huts, if they do, they may immediately explore those huts.'  # Zulu

# Widen the "Starting Tech" column to make room next to the new Description column
$ws.Columns.Item(2).ColumnWidth = 19.1821862348178

# Restore focus/selection to match the authored workbook state
$ws.Activate()
$ws.Range('C17').Select()
